$wb = $excel.ActiveWorkbook

# --- Insert new "Texas Notes" worksheet right after "About" ---
$aboutSheet = $wb.Worksheets.Item("About")
$notesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $aboutSheet)
$notesSheet.Name = "Texas Notes"

# --- Populate "Texas Notes" sheet content ---
$notesSheet.Range("A1").Value = "A more recent study from DNVGL "
$notesSheet.Range("A2").Value = "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html"
[void]$notesSheet.Hyperlinks.Add($notesSheet.Range("A2"), "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html")
$notesSheet.Range("A3").Value = "assumes a learning rate of 15-20% - closer to what we see in the renewables and storage industries."
$notesSheet.Range("A5").Value = "This suggests that the learning rate may be higher than the 2013 report used by EI. "
$notesSheet.Range("A6").Value = "So, we can take an average of some of these values just to be conservative"
$notesSheet.Range("B7").Value = "average learning rate"

$notesSheet.Range("A8").Value = "2018 DNVGL"
$notesSheet.Range("A8").HorizontalAlignment = -4131
$notesSheet.Range("B8").Formula = "=AVERAGE(0.15, 0.2)"

$notesSheet.Range("A9").Value = "2013 CRS report"
$notesSheet.Range("A9").HorizontalAlignment = -4131
$notesSheet.Range("B9").Value = 0.13

$notesSheet.Range("A10").Value = "average"
$notesSheet.Range("A10").HorizontalAlignment = -4152
$notesSheet.Range("B10").Formula = "=AVERAGE(B8:B9)"
$notesSheet.Range("B10").Interior.Color = 65535

$notesSheet.Columns("A").ColumnWidth = 16.8

# --- Point PDiCECpDoC's learning-rate cell at the new sheet's computed average ---
$pdiSheet = $wb.Worksheets.Item("PDiCECpDoC")
$pdiSheet.Range("B2").Formula = "='Texas Notes'!B10"

# --- Restore per-sheet selections / active tab ---
[void]$aboutSheet.Select()
[void]$aboutSheet.Range("E24").Select()

[void]$notesSheet.Select()
[void]$notesSheet.Range("G13").Select()

[void]$pdiSheet.Select()
[void]$pdiSheet.Range("C7").Select()

Write-Output "done"
